$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New record (row 84) appended by the automatic map update.
# Text-like columns are formatted as Text first so Excel stores them
# as strings (matching the source data) instead of auto-converting
# numeric-looking / date-looking text into numbers or dates.
$ws.Range("A84:H84").NumberFormat = "@"
$ws.Range("J84:L84").NumberFormat = "@"
$ws.Range("O84:P84").NumberFormat = "@"

$ws.Cells.Item(84, 1).Value = "7024"
$ws.Cells.Item(84, 2).Value = "8/25/2025"
$ws.Cells.Item(84, 3).Value = "SAAVEDRA 869"
$ws.Cells.Item(84, 4).Value = "3"
$ws.Cells.Item(84, 5).Value = "809155616"
$ws.Cells.Item(84, 6).Value = "PEBCOM"
$ws.Cells.Item(84, 7).Value = "Pendiente"
$ws.Cells.Item(84, 8).Value = "Picada"
$ws.Cells.Item(84, 9).Value = 1
$ws.Cells.Item(84, 10).Value = "Cambio"
$ws.Cells.Item(84, 11).Value = "Sin equipos"
$ws.Cells.Item(84, 12).Value = "Pasante"
$ws.Cells.Item(84, 13).Value = -58.402244
$ws.Cells.Item(84, 14).Value = -34.619401
$ws.Cells.Item(84, 15).Value = "San Telmo"
$ws.Cells.Item(84, 16).Value = "Capital Sur"
